# Scheduled runner: refresh market-price derived columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# across Leve-profit worksheets (ALC, BSM, CRP, CUL, GSM, LTW, WVR) with latest pulled data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 832.5
$ws.Range("I2").Value = 1076.6666
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 1076.6666
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = -963.6666
$ws.Range("N2").Value = -326

$ws.Range("H111").Value = 2572.5715
$ws.Range("I111").Value = 2501.3333
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 7503.999899999999
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = -4436.999899999999
$ws.Range("N111").Value = -15134

$ws.Range("H118").Value = 7120
$ws.Range("I118").Value = 906.6667
$ws.Range("J118").Value = 13333.333
$ws.Range("K118").Value = 2720.0001
$ws.Range("L118").Value = 39999.999
$ws.Range("M118").Value = -1063.0001
$ws.Range("N118").Value = -43313.999

$ws.Range("H125").Value = 2400
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 3000
$ws.Range("K125").Value = 18000
$ws.Range("L125").Value = 27000
$ws.Range("M125").Value = -15540
$ws.Range("N125").Value = -31920

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 117.71429
$ws.Range("J80").Value = 121
$ws.Range("L80").Value = 121
$ws.Range("N80").Value = -2117

$ws.Range("H83").Value = 117.71429
$ws.Range("J83").Value = 121
$ws.Range("L83").Value = 605
$ws.Range("N83").Value = -10589

$ws.Range("H86").Value = 100001250
$ws.Range("I86").Value = 200000000
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 200000000
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -199998877
$ws.Range("N86").Value = -4746

$ws.Range("H89").Value = 100001250
$ws.Range("I89").Value = 200000000
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 1000000000
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -999994384
$ws.Range("N89").Value = -23732

$ws.Range("H94").Value = 13911.6875
$ws.Range("I94").Value = 1510.4615
$ws.Range("J94").Value = 67650.336
$ws.Range("K94").Value = 1510.4615
$ws.Range("L94").Value = 67650.336
$ws.Range("M94").Value = -1059.4615
$ws.Range("N94").Value = -68552.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3237
$ws.Range("I62").Value = 3002.5
$ws.Range("J62").Value = 3393.3333
$ws.Range("K62").Value = 3002.5
$ws.Range("L62").Value = 3393.3333
$ws.Range("M62").Value = -2378.5
$ws.Range("N62").Value = -4641.3333

$ws.Range("H65").Value = 3237
$ws.Range("I65").Value = 3002.5
$ws.Range("J65").Value = 3393.3333
$ws.Range("K65").Value = 15012.5
$ws.Range("L65").Value = 16966.6665
$ws.Range("M65").Value = -11892.5
$ws.Range("N65").Value = -23206.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1088.8276
$ws.Range("I107").Value = 295.16666
$ws.Range("J107").Value = 1295.8695
$ws.Range("K107").Value = 885.4999799999999
$ws.Range("L107").Value = 3887.6085
$ws.Range("M107").Value = 1034.50002
$ws.Range("N107").Value = -7727.6085

$ws.Range("H113").Value = 1623.0834
$ws.Range("I113").Value = 797.25
$ws.Range("J113").Value = 2036
$ws.Range("K113").Value = 2391.75
$ws.Range("L113").Value = 6108
$ws.Range("M113").Value = -221.75
$ws.Range("N113").Value = -10448

$ws.Range("H131").Value = 885.0933
$ws.Range("J131").Value = 979.6070999999999
$ws.Range("L131").Value = 2938.8213
$ws.Range("N131").Value = -13018.8213

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

$ws.Range("H25").Value = 3999
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 3999
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 3999
$ws.Range("N25").Value = -5057
$ws.Range("M25").ClearContents()

$ws.Range("H70").Value = 25005002
$ws.Range("I70").Value = 44449276
$ws.Range("J70").Value = 5219.2856
$ws.Range("K70").Value = 44449276
$ws.Range("L70").Value = 5219.2856
$ws.Range("M70").Value = -44449006
$ws.Range("N70").Value = -5759.2856

$ws.Range("H73").Value = 25005002
$ws.Range("I73").Value = 44449276
$ws.Range("J73").Value = 5219.2856
$ws.Range("K73").Value = 44449276
$ws.Range("L73").Value = 5219.2856
$ws.Range("M73").Value = -44448340
$ws.Range("N73").Value = -7091.2856

$ws.Range("H80").Value = 2773.8235
$ws.Range("I80").Value = 2767.5
$ws.Range("J80").Value = 2777.2727
$ws.Range("K80").Value = 2767.5
$ws.Range("L80").Value = 2777.2727
$ws.Range("M80").Value = -1769.5
$ws.Range("N80").Value = -4773.2727

$ws.Range("H83").Value = 2773.8235
$ws.Range("I83").Value = 2767.5
$ws.Range("J83").Value = 2777.2727
$ws.Range("K83").Value = 13837.5
$ws.Range("L83").Value = 13886.3635
$ws.Range("M83").Value = -8845.5
$ws.Range("N83").Value = -23870.3635

$ws.Range("H107").Value = 627.8570999999999
$ws.Range("I107").Value = 905.5
$ws.Range("J107").Value = 457
$ws.Range("K107").Value = 905.5
$ws.Range("L107").Value = 457
$ws.Range("M107").Value = 1014.5
$ws.Range("N107").Value = -4297

$ws.Range("H122").Value = 1025
$ws.Range("J122").Value = 1025
$ws.Range("L122").Value = 3075
$ws.Range("N122").Value = -7975

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3972.7917
$ws.Range("I40").Value = 4146.8125
$ws.Range("J40").Value = 3624.75
$ws.Range("K40").Value = 4146.8125
$ws.Range("L40").Value = 3624.75
$ws.Range("M40").Value = -4010.8125
$ws.Range("N40").Value = -3896.75

$ws.Range("H93").Value = 1357.9474
$ws.Range("I93").Value = 1338.2307
$ws.Range("J93").Value = 1400.6666
$ws.Range("K93").Value = 1338.2307
$ws.Range("L93").Value = 1400.6666
$ws.Range("M93").Value = -90.23070000000007
$ws.Range("N93").Value = -3896.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 5001500
$ws.Range("J2").Value = 5001500
$ws.Range("L2").Value = 5001500
$ws.Range("N2").Value = -5001724

Write-Host "Updated market price data across worksheets."
